{"js": "// Insert a new list item \"Add the Division method and DivisionTestCase5.\"\n// right after the \"Right-click, then select Run Tests...\" paragraph\n// (and before the trailing empty list paragraph), matching the\n// existing ListParagraph style + numbering (ilvl 0, numId 1).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the anchor paragraph by its text.\nconst anchorText = \"Right-click, then select Run Tests. This will open the Test Explorer.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// insertParagraph(\"After\") clones the anchor paragraph's formatting\n// (style + numbering), exactly like pressing Enter at the end of it.\nanchor.insertParagraph(\"Add the Division method and DivisionTestCase5.\", \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new list item \"Add the Division method and DivisionTestCase5.\"\n# right after the \"Right-click, then select Run Tests...\" paragraph\n# (and before the trailing empty list paragraph), matching the\n# existing ListParagraph style + numbering (ilvl 0, numId 1).\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph using Find.\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute(\"Right-click, then select Run Tests. This will open the Test Explorer.\")\nif (-not $found) {\n    throw \"Could not find anchor paragraph text.\"\n}\n$anchorStart = $findRange.Start\n\n# Resolve the Find hit to its containing Paragraphs() index so we can\n# reliably insert/set text on real paragraph objects (Range objects\n# returned by Find can behave oddly once collapsed/mutated).\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Start -le $anchorStart -and $anchorStart -lt $p.Range.End) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not resolve anchor paragraph index.\"\n}\n\n$anchorPara = $d.Paragraphs.Item($anchorIndex)\n\n# InsertParagraphAfter clones the anchor paragraph's formatting\n# (style + numbering), exactly like pressing Enter at the end of it.\n$anchorPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($anchorIndex + 1)\n$newPara.Range.Text = \"Add the Division method and DivisionTestCase5.\"\n"}
